$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold / bordered / centered-top style on B1 first
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160

# Copy that exact formatting onto A2 (keeps a single shared style index,
# instead of rebuilding it property-by-property which forks a new style)
$r2 = $ws.Range("A2")
$r1.Copy()
$r2.PasteSpecial(-4122)
